$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Update the repayment strategy value (B17) from "RBI (India)" to
# "Overdue/Due Fee/Int,Principal"
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Give B17 a distinct look (green fill, Arial 10) to flag the new scenario
$ws.Range("B17").Interior.Color = 5296274
$ws.Range("B17").Font.Name = "Arial"
$ws.Range("B17").Font.Size = 10

# Move the active selection to B17, matching the edited workbook
$ws.Activate()
$ws.Range("B17").Select()
